# Lineup editor: append newly-scraped bands to the bottom of the list.
# (Only the band name is known for these entries; the remaining columns
# - country, genre, youtube url, votes, about - are left blank for now.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bands = @(
    "LEPROUS",
    "ALESTORM",
    "ELUVEITIE",
    "BEHEMOTH",
    "CORONER",
    "ENSIFERUM",
    "ACCEPT",
    "WATAIN",
    "LOUDNESS",
    "IGORRR",
    "KATAKLYSM",
    "M2TM",
    "OBITUARY",
    "BLACK STAR RIDERS",
    "HATEBREED",
    "JUDAS PRIEST",
    "MUNICIPAL WASTE",
    "CANNIBAL CORPSE",
    "EPICA",
    "CHILDREN OF BODOM"
)

$startRow = 118
$row = $startRow
foreach ($band in $bands) {
    $ws.Cells.Item($row, 1).Value = $band
    $row = $row + 1
}

# Scroll the sheet down to the freshly-added rows and leave the cursor on
# the first new row, column D (mirrors where editing left off).
$ws.Range("D" + $startRow).Select()
$excel.ActiveWindow.ScrollRow = 106
